$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the literal string (incl. trailing zeros / grouping) is kept.
$priceTextForced = [ordered]@{
    "D4" = '1.005'
    "D5" = '313.42'
    "D6" = '1.004'
    "D7" = '0.3877'
    "D8" = '0.4011'
    "D9" = '1.485'
    "D10" = '1.004'
    "D11" = '52.35'
    "D12" = '0.08737'
    "D13" = '24.97'
    "D14" = '7.478'
    "D15" = '8.017'
    "D16" = '0.00001342'
    "D18" = '97.94'
    "D19" = '0.07077'
    "D20" = '19.91'
    "D21" = '7.220'
    "D22" = '1.004'
    "D23" = '14.19'
    "D25" = '2.354'
    "D26" = '2.930'
    "D27" = '22.56'
    "D28" = '162.62'
    "D29" = '8.417'
    "D30" = '136.51'
    "D31" = '5.192'
    "D33" = '0.08779'
    "D34" = '7.357'
    "D35" = '1.020'
    "D36" = '0.2795'
    "D37" = '1.947'
    "D38" = '0.02882'
    "D39" = '10.69'
    "D41" = '0.09091'
    "D42" = '0.7855'
    "D43" = '1.451'
    "D44" = '16.64'
    "D45" = '0.7221'
    "D46" = '2.582'
    "D47" = '4.194'
    "D48" = '1.003'
    "D49" = '1.339'
    "D50" = '138.20'
    "D51" = '0.08007'
}
foreach ($addr in $priceTextForced.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceTextForced[$addr]
}

# Cells whose new text is already unambiguous (e.g. contains two "." separators)
# and stays text on assignment without any format change.
$priceTextNatural = [ordered]@{
    "D2" = '24.457.57'
    "D3" = '1.687.93'
    "D17" = '1.685.41'
    "D24" = '24.436.66'
    "D32" = '1.868.50'
}
foreach ($addr in $priceTextNatural.Keys) {
    $ws.Range($addr).Value = $priceTextNatural[$addr]
}

# --- Column E (Volume(1h)) updates ---
# Values are padded with spaces (e.g. "  -1.07%  ") so they are never number-like.
$volumeUpdates = [ordered]@{
    "E2" = '  -1.07%  '
    "E3" = '  -0.52%  '
    "E4" = '  +0.10%  '
    "E5" = '  -1.32%  '
    "E6" = '  +0.10%  '
    "E7" = '  -2.04%  '
    "E8" = '  -1.71%  '
    "E9" = '  -0.33%  '
    "E10" = '  +0.05%  '
    "E11" = '  -0.88%  '
    "E12" = '  -2.04%  '
    "E13" = '  +6.07%  '
    "E14" = '  +2.93%  '
    "E15" = '  -0.17%  '
    "E16" = '  +1.40%  '
    "E17" = '  -0.74%  '
    "E18" = '  -2.12%  '
    "E19" = '  +0.47%  '
    "E20" = '  +1.23%  '
    "E21" = '  +3.08%  '
    "E22" = '  +0.18%  '
    "E23" = '  -1.08%  '
    "E24" = '  -1.08%  '
    "E25" = '  -0.51%  '
    "E26" = '  -10.45%  '
    "E27" = '  -0.76%  '
    "E28" = '  +0.17%  '
    "E29" = '  +11.68%  '
    "E30" = '  -0.18%  '
    "E31" = '  +0.55%  '
    "E32" = '  -0.86%  '
    "E33" = '  +0.78%  '
    "E34" = '  +3.68%  '
    "E35" = '  -3.47%  '
    "E36" = '  +1.87%  '
    "E37" = '  +3.39%  '
    "E38" = '  +5.60%  '
    "E39" = '  -6.66%  '
    "E40" = '  -2.35%  '
    "E41" = '  -1.60%  '
    "E42" = '  +2.45%  '
    "E43" = '  -1.77%  '
    "E44" = '  +3.13%  '
    "E45" = '  +0.56%  '
    "E46" = '  -0.42%  '
    "E47" = '  -0.74%  '
    "E48" = '  +0.01%  '
    "E49" = '  +1.44%  '
    "E50" = '  -1.55%  '
    "E51" = '  +0.28%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}